$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '23.703.84'
$ws.Range("E2").Value = '  +1.21%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.655.83'
$ws.Range("E3").Value = '  +1.30%  '

$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("E5").Value = '  +0.06%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '303.16'
$ws.Range("E6").Value = '  -0.07%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3805'
$ws.Range("E7").Value = '  +0.62%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3622'
$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '51.27'
$ws.Range("E9").Value = '  -0.63%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.247'
$ws.Range("E10").Value = '  +1.72%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08220'
$ws.Range("E11").Value = '  +0.56%  '

$ws.Range("E12").Value = '  +0.06%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.67'
$ws.Range("E13").Value = '  +1.47%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.528'
$ws.Range("E14").Value = '  +1.04%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.443'
$ws.Range("E15").Value = '  +0.92%  '

$ws.Range("E16").Value = '  -0.27%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.656.97'
$ws.Range("E17").Value = '  +1.71%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '97.47'
$ws.Range("E18").Value = '  +2.53%  '

$ws.Range("E19").Value = '  +1.14%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.814'
$ws.Range("E20").Value = '  +3.52%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.71'
$ws.Range("E21").Value = '  +1.52%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.002'
$ws.Range("E22").Value = '  +0.14%  '

$ws.Range("E23").Value = '  +2.61%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '23.713.96'
$ws.Range("E24").Value = '  +1.26%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.517'
$ws.Range("E25").Value = '  +0.01%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.058'
$ws.Range("E26").Value = '  +0.38%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.29'
$ws.Range("E27").Value = '  +0.90%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '153.63'
$ws.Range("E28").Value = '  +1.79%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.229'

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '134.85'
$ws.Range("E30").Value = '  +1.18%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.841.34'
$ws.Range("E31").Value = '  +1.63%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.957'
$ws.Range("E32").Value = '  +5.46%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.217'
$ws.Range("E33").Value = '  +2.35%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.068'
$ws.Range("E34").Value = '  +2.19%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '11.76'
$ws.Range("E35").Value = '  +4.87%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02828'
$ws.Range("E36").Value = '  +3.09%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2535'
$ws.Range("E37").Value = '  +1.85%  '

$ws.Range("B38").Value = 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.116'
$ws.Range("E38").Value = '  +1.81%  '

$ws.Range("B39").Value = 'Stellar'
$ws.Range("C39").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.08796'
$ws.Range("E39").Value = '  +0.34%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.07123'
$ws.Range("E40").Value = '  +0.30%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '13.00'
$ws.Range("E41").Value = '  +7.11%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7049'
$ws.Range("E42").Value = '  +0.74%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.335'
$ws.Range("E43").Value = '  -0.04%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.09'
$ws.Range("E44").Value = '  +2.17%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6527'
$ws.Range("E45").Value = '  +0.55%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.320'
$ws.Range("E46").Value = '  +2.20%  '

$ws.Range("E47").Value = '  -0.01%  '

$ws.Range("E48").Value = '  +0.29%  '

$ws.Range("E49").Value = '  -0.15%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '128.27'
$ws.Range("E50").Value = '  +0.93%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.190'
$ws.Range("E51").Value = '  +0.09%  '
